# This workbook tracks survey-style preference data on 7 sheets (DIP, FP,
# SP, TP, FoP, FiP, SiP). Each sheet has the same layout: rows 2-12 hold one
# "response count" record each in column A, with preference-percentage
# inputs in J:O and derived formulas in B:H (per-choice allocation) and P
# (checksum). The edit:
#   1. Removes the records whose "No. of responses" (col A) were 93, 103
#      and 108 (the old rows 8, 10 and 11) from every sheet.
#   2. Appends a brand-new record (A = 133) with its own J:O inputs and the
#      same B:H / P formulas, landing in row 10 after the deletions above.
#   3. Leaves the final selection wherever the author's last click landed on
#      each sheet (cosmetic, but reproduced for fidelity).

$wb = $excel.ActiveWorkbook

# sheet index (1-based, left-to-right tab order) -> new record's J:O inputs.
# Column A for the new record is always 133 on every sheet.
$newRows = @{
    1 = @{ J = 14.3;  K = 12;   L = 34.6;               M = 15.8;               N = 10.5;               O = 11.3 }
    2 = @{ J = 14.3;  K = 4.5;  L = 62.4;                M = 0;                  N = 16.5;               O = 2.2999999999999998 }
    3 = @{ J = 12;    K = 11.3; L = 18;                  M = 0;                  N = 52.6;               O = 6 }
    4 = @{ J = 34.6;  K = 16.5; L = 10.5;                M = 13.5;               N = 18;                 O = 6.8 }
    5 = @{ J = 15.8;  K = 26.3; L = 5.3;                 M = 33.799999999999997; N = 9;                  O = 9.8000000000000007 }
    6 = @{ J = 10.5;  K = 26.3; L = 4.5;                 M = 24.1;               N = 2.2999999999999998; O = 32.200000000000003 }
    7 = @{ J = 11.3;  K = 15.8; L = 1.5;                 M = 26.3;               N = 2.2999999999999998; O = 42.9 }
}

# Final post-edit selection per sheet (cosmetic, mirrors the author's state).
$selections = @{
    1 = "B17"
    2 = "A9:XFD10"
    3 = "B19"
    4 = "A8:XFD8"
    5 = "A9:XFD10"
    6 = "A9:XFD10"
    7 = "A9:XFD10"
}

for ($i = 1; $i -le 7; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Delete the three obsolete records (old rows 11, 10, 8 — highest first
    # so earlier deletions don't shift the still-to-be-deleted row numbers).
    $ws.Rows.Item(11).Delete()
    $ws.Rows.Item(10).Delete()
    $ws.Rows.Item(8).Delete()

    # After the three deletions, rows 8 and 9 hold the old 99- and
    # 110-response records (shifted up) and row 10 is free for the new
    # 133-response record.
    $data = $newRows[$i]
    $ws.Range("A10").Value = 133
    $ws.Range("J10").Value = $data.J
    $ws.Range("K10").Value = $data.K
    $ws.Range("L10").Value = $data.L
    $ws.Range("M10").Value = $data.M
    $ws.Range("N10").Value = $data.N
    $ws.Range("O10").Value = $data.O

    $ws.Range("B10").Formula = "=A10*J10%"
    $ws.Range("C10").Formula = "=A10*K10%"
    $ws.Range("D10").Formula = "=A10*L10%"
    $ws.Range("E10").Formula = "=A10*M10%"
    $ws.Range("F10").Formula = "=A10*N10%"
    $ws.Range("G10").Formula = "=A10*O10%"
    $ws.Range("H10").Formula = "=SUM(B10:G10)"
    $ws.Range("P10").Formula = "=SUM(J10:O10)"

    # Match the number format (-> style) already used by the rows above so
    # the new row renders the same way (integer-style numeric format).
    $ws.Range("B10:H10").NumberFormat = "0"

    # Restore the sheet's final on-screen selection.
    $ws.Range($selections[$i]).Select()
}

$dip = $wb.Worksheets.Item(1)
$dip.Activate()
